# Update the "dSF" (F) column values for several rows.
# These edits correspond to a repull/recalculation of the data
# (mean calculation applied to the F column), while the E column
# ("dS0") values stay untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 4
    3  = -2
    8  = -1
    13 = -5
    14 = -16
    15 = -2
    16 = -3
    17 = 1
    20 = 0
    22 = -3
    23 = -5
    25 = -1
    27 = -1
    28 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
